# Generate Report for Handoff
# The "0034157f-..." file is being handed off again, so its "Latest Handoff
# Datetime" (column D) is refreshed with a new timestamp on both the
# zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-03-04 10:09:27"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-03-04 10:09:39"
